$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.709.57"
$ws.Range("E2").Value = "  -3.29%  "

$ws.Range("D3").Value = "3.169.66"
$ws.Range("E3").Value = "  -7.87%  "

$ws.Range("E4").Value = "  +0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "564.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.84%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "170.86"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.89%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.58%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "3.167.40"
$ws.Range("E9").Value = "  -7.82%  "

$ws.Range("E10").Value = "  -6.53%  "

$ws.Range("E11").Value = "  -5.72%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.394"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.79%  "

$ws.Range("D13").Value = "3.725.90"
$ws.Range("E13").Value = "  -7.67%  "

$ws.Range("E14").Value = "  +0.58%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.41"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -8.86%  "

$ws.Range("D16").Value = "64.602.16"
$ws.Range("E16").Value = "  -3.33%  "

$ws.Range("E17").Value = "  -6.79%  "

$ws.Range("D18").Value = "3.163.39"
$ws.Range("E18").Value = "  -8.44%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.96%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.00"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -6.48%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "355.78"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.63%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.25"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.57%  "

$ws.Range("E23").Value = "  +1.07%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "68.62"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -6.51%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.501"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.90%  "

$ws.Range("E26").Value = "  -9.58%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.65"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.69%  "

$ws.Range("E28").Value = "  -2.14%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("E30").Value = "  -0.15%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.57%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -8.04%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "22.02"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -7.09%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.68"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.31%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -6.07%  "

$ws.Range("E36").Value = "  -8.65%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "153.50"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.77%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.829"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.91%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "26.46"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.66%  "

$ws.Range("E40").Value = "  -5.10%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.52"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.75%  "

$ws.Range("D42").Value = "2.651.00"
$ws.Range("E42").Value = "  -3.60%  "

$ws.Range("E43").Value = "  -7.49%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "6.03"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.40%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "39.37"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.04%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "24.22"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -6.15%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0654"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.61%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "321.85"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.20%  "

$ws.Range("E49").Value = "  -5.14%  "

$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("E51").Value = "  +0.02%  "
